$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4092030496915097
$ws.Range("D2").Value = 0.6863471501900005

$ws.Range("C3").Value = -1.186590314980815
$ws.Range("D3").Value = 0.2480455526801508

$ws.Range("C4").Value = -0.9545872097143744
$ws.Range("D4").Value = 0.3501586132443055

$ws.Range("C5").Value = -0.5560064002629599
$ws.Range("D5").Value = 0.583817679544365

$ws.Range("C6").Value = -1.037338867522083
$ws.Range("D6").Value = 0.3108458365873463

$ws.Range("C7").Value = -0.8556429826365585
$ws.Range("D7").Value = 0.4014184801490595

$ws.Range("C8").Value = -0.7428556367503022
$ws.Range("D8").Value = 0.4654280414399126

$ws.Range("C9").Value = 0.09830806759517204
$ws.Range("D9").Value = 0.9225778896668391

$ws.Range("C10").Value = 0.252114352602187
$ws.Range("D10").Value = 0.8032930674845249

$ws.Range("C11").Value = 0.1432549447867733
$ws.Range("D11").Value = 0.8873928631734185
